$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is text that Excel would NOT misinterpret as a number
# (multiple dots, subscript digits, etc.) -> can be set directly as .Value
$textCells = @{
    'D2' = '69.618.38'
    'E2' = '  +0.27%  '
    'D3' = '2.501.65'
    'E3' = '  -0.08%  '
    'E4' = '  +0.02%  '
    'E5' = '  -0.27%  '
    'E6' = '  +0.32%  '
    'E7' = '  -0.01%  '
    'E8' = '  -1.58%  '
    'D9' = '2.500.67'
    'E9' = '  -0.07%  '
    'E10' = '  +1.00%  '
    'E11' = '  +0.19%  '
    'E12' = '  +3.06%  '
    'E13' = '  +1.49%  '
    'D14' = '2.959.45'
    'E14' = '  -0.09%  '
    'D15' = '69.490.89'
    'E15' = '  +0.16%  '
    'E16' = '  +1.92%  '
    'E17' = '  -0.63%  '
    'D18' = '2.493.00'
    'E18' = '  -0.48%  '
    'E19' = '  -0.96%  '
    'E20' = '  -4.55%  '
    'E21' = '  -0.03%  '
    'E22' = '  -0.86%  '
    'E23' = '  -0.16%  '
    'E24' = '  -0.11%  '
    'E25' = '  +2.74%  '
    'E26' = '  -0.57%  '
    'E27' = '  -1.55%  '
    'D28' = '2.627.21'
    'E28' = '  -0.20%  '
    'E29' = '  +0.14%  '
    'D30' = '0.0₃0892'
    'E30' = '  -1.12%  '
    'E31' = '  -1.00%  '
    'E32' = '  -0.93%  '
    'E33' = '  -3.25%  '
    'E34' = '  -1.05%  '
    'E35' = '  -0.02%  '
    'E36' = '  -1.18%  '
    'E37' = '  +2.37%  '
    'E38' = '  +0.63%  '
    'E39' = '  +0.32%  '
    'E40' = '  +0.03%  '
    'E41' = '  -0.11%  '
    'E43' = '  +0.05%  '
    'E44' = '  +0.19%  '
    'E45' = '  -3.93%  '
    'E46' = '  -6.15%  '
    'E47' = '  -1.22%  '
    'E48' = '  -0.44%  '
    'E49' = '  -1.70%  '
    'E50' = '  +0.24%  '
    'E51' = '  -0.71%  '
}
foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Value = $textCells[$addr]
}

# Cells whose new value looks like a plain number (e.g. "576.08") but must
# remain stored as TEXT, matching the original inline-string cell content.
# Setting NumberFormat to "@" (Text) before assigning the value forces Excel
# to keep it as a string instead of silently converting it to a numeric cell;
# resetting the Style back to "Normal" afterwards removes the now-unneeded
# explicit number-format style so the cell keeps its original (default) style.
$numericLookingTextCells = @{
    'D5' = '576.08'
    'D6' = '166.74'
    'D13' = '4.94'
    'D17' = '24.72'
    'D19' = '11.21'
    'D20' = '7.48'
    'D21' = '348.26'
    'D23' = '1.93'
    'D25' = '70.47'
    'D29' = '0.999'
    'D32' = '457.79'
    'D33' = '1.21'
    'D36' = '0.117'
    'D37' = '156.79'
    'D38' = '19.06'
    'D44' = '38.17'
    'D45' = '2.20'
    'D47' = '141.51'
    'D50' = '0.0733'
}
foreach ($addr in $numericLookingTextCells.Keys) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $numericLookingTextCells[$addr]
    $c.Style = "Normal"
}
